$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Fix the typo spanning the "ФИО" / "новер" run boundary.
#    Original combined text (run-split in the middle of "ФИО"):
#      "...содержащую ФИ" + "О студента, номер группы и " + "новер" + " варианта..."
#    Target combined text:
#      "...содержащую ФИ" + "О студента, номер группы и номер" + " варианта..."
#    We match across the existing run/proofErr boundary so the edit
#    naturally swallows the <w:proofErr w:type="spellStart"/> marker
#    that wrapped "новер".
# ------------------------------------------------------------------
$full = $d.Content
$full.Find.Execute("О студента, номер группы и новер", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPos = $full.Start
$full.Text = "О студента, номер группы и номер"
$mergedEnd = $full.End

# ------------------------------------------------------------------
# 2. Remove the now-orphaned <w:proofErr w:type="spellEnd"/> that used
#    to close out the "новер" spell-check span (it now sits right
#    after $mergedEnd, before " варианта..."). Round-tripping a tiny
#    range that straddles that boundary through a throw-away edit
#    forces the run/proofErr markup there to be rebuilt without it.
# ------------------------------------------------------------------
$tailProbe = $d.Range($mergedEnd - 1, $mergedEnd + 1)
$tailText = $tailProbe.Text
$tailProbe.Text = "zz"
$tailRestore = $d.Range($mergedEnd - 1, $mergedEnd - 1 + $tailText.Length)
$tailRestore.Text = $tailText

# ------------------------------------------------------------------
# 3. Re-introduce the run break between "...содержащую ФИ" and
#    "О студента, номер группы и номер" (the editor's cursor sat there)
#    by bracketing it with a throw-away bookmark, then drop a fresh
#    "_GoBack" bookmark at the edit's end point (right before
#    " варианта..."), matching Word's own "last edit" bookmark.
# ------------------------------------------------------------------
$splitRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("zzTempSplit", $splitRange)

$goBackRange = $d.Range($mergedEnd, $mergedEnd)
$d.Bookmarks.Add("_GoBack", $goBackRange)

$d.Bookmarks("zzTempSplit").Delete()
